$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 4).Value = "'64.125.25"
$ws.Cells.Item(2, 5).Value = "'  +1.22%  "
$ws.Cells.Item(3, 4).Value = "'2.780.55"
$ws.Cells.Item(3, 5).Value = "'  +2.11%  "
$ws.Cells.Item(4, 5).Value = "'  -0.35%  "
$ws.Cells.Item(5, 4).Value = "'589.47"
$ws.Cells.Item(5, 5).Value = "'  +0.68%  "
$ws.Cells.Item(6, 4).Value = "'160.63"
$ws.Cells.Item(6, 5).Value = "'  +6.36%  "
$ws.Cells.Item(7, 4).Value = "'0.622"
$ws.Cells.Item(7, 5).Value = "'  +2.24%  "
$ws.Cells.Item(8, 5).Value = "'  +0.34%  "
$ws.Cells.Item(9, 5).Value = "'  +0.81%  "
$ws.Cells.Item(10, 5).Value = "'  +1.05%  "
$ws.Cells.Item(11, 4).Value = "'0.397"
$ws.Cells.Item(11, 5).Value = "'  +1.84%  "
$ws.Cells.Item(12, 5).Value = "'  +1.12%  "
$ws.Cells.Item(13, 4).Value = "'3.277.64"
$ws.Cells.Item(13, 5).Value = "'  +2.14%  "
$ws.Cells.Item(14, 4).Value = "'27.49"
$ws.Cells.Item(14, 5).Value = "'  +2.39%  "
$ws.Cells.Item(15, 4).Value = "'64.039.87"
$ws.Cells.Item(15, 5).Value = "'  +1.22%  "
$ws.Cells.Item(16, 4).Value = "'0.0000160"
$ws.Cells.Item(16, 5).Value = "'  +5.49%  "
$ws.Cells.Item(17, 4).Value = "'2.788.08"
$ws.Cells.Item(17, 5).Value = "'  +1.68%  "
$ws.Cells.Item(18, 4).Value = "'12.46"
$ws.Cells.Item(18, 5).Value = "'  +4.10%  "
$ws.Cells.Item(19, 5).Value = "'  +3.68%  "
$ws.Cells.Item(20, 4).Value = "'367.46"
$ws.Cells.Item(20, 5).Value = "'  +0.83%  "
$ws.Cells.Item(21, 4).Value = "'7.05"
$ws.Cells.Item(21, 5).Value = "'  +0.17%  "
$ws.Cells.Item(22, 4).Value = "'0.576"
$ws.Cells.Item(22, 5).Value = "'  +7.51%  "
$ws.Cells.Item(23, 5).Value = "'  -0.18%  "
$ws.Cells.Item(24, 4).Value = "'67.45"
$ws.Cells.Item(24, 5).Value = "'  +2.70%  "
$ws.Cells.Item(25, 4).Value = "'0.177"
$ws.Cells.Item(25, 5).Value = "'  +6.01%  "
$ws.Cells.Item(26, 4).Value = "'8.85"
$ws.Cells.Item(26, 5).Value = "'  +3.46%  "
$ws.Cells.Item(27, 4).Value = "'0.0₃0971"
$ws.Cells.Item(27, 5).Value = "'  +11.70%  "
$ws.Cells.Item(28, 5).Value = "'  +0.24%  "
$ws.Cells.Item(29, 4).Value = "'2.07"
$ws.Cells.Item(29, 5).Value = "'  +1.58%  "
$ws.Cells.Item(30, 5).Value = "'  +1.80%  "
$ws.Cells.Item(31, 4).Value = "'1.27"
$ws.Cells.Item(31, 5).Value = "'  +6.47%  "
$ws.Cells.Item(32, 4).Value = "'5.18"
$ws.Cells.Item(32, 5).Value = "'  +8.47%  "
$ws.Cells.Item(33, 4).Value = "'170.25"
$ws.Cells.Item(33, 5).Value = "'  -0.47%  "
$ws.Cells.Item(34, 4).Value = "'20.91"
$ws.Cells.Item(34, 5).Value = "'  +1.55%  "
$ws.Cells.Item(35, 5).Value = "'  -0.03%  "
$ws.Cells.Item(36, 5).Value = "'  +4.36%  "
$ws.Cells.Item(37, 5).Value = "'  +1.93%  "
$ws.Cells.Item(38, 5).Value = "'  +1.20%  "
$ws.Cells.Item(39, 2).Value = "RenderToken"
$ws.Cells.Item(39, 3).Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Cells.Item(39, 4).Value = "'6.32"
$ws.Cells.Item(39, 5).Value = "'  +11.56%  "
$ws.Cells.Item(40, 2).Value = "Filecoin"
$ws.Cells.Item(40, 3).Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Cells.Item(40, 4).Value = "'4.29"
$ws.Cells.Item(40, 5).Value = "'  +0.44%  "
$ws.Cells.Item(41, 2).Value = "Bittensor"
$ws.Cells.Item(41, 3).Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Cells.Item(41, 4).Value = "'342.34"
$ws.Cells.Item(41, 5).Value = "'  -2.64%  "
$ws.Cells.Item(42, 4).Value = "'40.27"
$ws.Cells.Item(42, 5).Value = "'  +2.46%  "
$ws.Cells.Item(43, 4).Value = "'22.50"
$ws.Cells.Item(43, 5).Value = "'  +0.86%  "
$ws.Cells.Item(44, 2).Value = "Hedera"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Cells.Item(44, 4).Value = "'0.0615"
$ws.Cells.Item(44, 5).Value = "'  +3.46%  "
$ws.Cells.Item(45, 2).Value = "InjectiveProtocol"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Cells.Item(45, 4).Value = "'22.53"
$ws.Cells.Item(45, 5).Value = "'  +2.81%  "
$ws.Cells.Item(46, 5).Value = "'  +1.58%  "
$ws.Cells.Item(47, 4).Value = "'0.0263"
$ws.Cells.Item(47, 5).Value = "'  +1.49%  "
$ws.Cells.Item(48, 4).Value = "'138.98"
$ws.Cells.Item(48, 5).Value = "'  -0.29%  "
$ws.Cells.Item(49, 4).Value = "'0.104"
$ws.Cells.Item(49, 5).Value = "'  +2.15%  "
$ws.Cells.Item(50, 4).Value = "'2.171.96"
$ws.Cells.Item(50, 5).Value = "'  +0.38%  "
$ws.Cells.Item(51, 5).Value = "'  +0.47%  "
